$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 122, pushing existing rows 122-133 down to 124-135.
$ws.Rows.Item(122).Resize(2).Insert()

# New row 122: Haba, "Primera", Provincia de Limarí, date 44461
$ws.Cells.Item(122,1).Value = 9
$ws.Cells.Item(122,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(122,3).Value = "Metropolitana"
$ws.Cells.Item(122,4).Value = 44461
$ws.Cells.Item(122,5).Value = 13
$ws.Cells.Item(122,6).Value = 100112026
$ws.Cells.Item(122,7).Value = "Haba"
$ws.Cells.Item(122,8).Value = "Sin especificar"
$ws.Cells.Item(122,9).Value = "Primera"
$ws.Cells.Item(122,10).Value = 52
$ws.Cells.Item(122,11).Value = 14000
$ws.Cells.Item(122,12).Value = 15000
$ws.Cells.Item(122,13).Value = 14500
$ws.Cells.Item(122,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(122,15).Value = "Provincia de Limarí"
$ws.Cells.Item(122,16).Value = 580
$ws.Cells.Item(122,17).Value = 25
$ws.Cells.Item(122,18).Value = "Hortaliza"

# New row 123: Haba, "Segunda", Provincia de Limarí, date 44461
$ws.Cells.Item(123,1).Value = 9
$ws.Cells.Item(123,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(123,3).Value = "Metropolitana"
$ws.Cells.Item(123,4).Value = 44461
$ws.Cells.Item(123,5).Value = 13
$ws.Cells.Item(123,6).Value = 100112026
$ws.Cells.Item(123,7).Value = "Haba"
$ws.Cells.Item(123,8).Value = "Sin especificar"
$ws.Cells.Item(123,9).Value = "Segunda"
$ws.Cells.Item(123,10).Value = 16
$ws.Cells.Item(123,11).Value = 11000
$ws.Cells.Item(123,12).Value = 12000
$ws.Cells.Item(123,13).Value = 11500
$ws.Cells.Item(123,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(123,15).Value = "Provincia de Limarí"
$ws.Cells.Item(123,16).Value = 460
$ws.Cells.Item(123,17).Value = 25
$ws.Cells.Item(123,18).Value = "Hortaliza"

Write-Output "rows inserted and populated"
